$wb = $excel.ActiveWorkbook

# --- Sheet: mmWave(InBed) -> add rows 76-88 ---
$ws10 = $wb.Worksheets.Item("mmWave(InBed)")
$ws10.Range("A76:A88").NumberFormat = "@"
$ws10.Cells.Item(76, 1).Value = "2026-02-01"
$ws10.Cells.Item(76, 2).Value = "20:18:35"
$ws10.Cells.Item(76, 3).Value = "20:00"
$ws10.Cells.Item(76, 4).Value = "Bedroom"
$ws10.Cells.Item(76, 5).Value = "In Bed"
$ws10.Cells.Item(76, 6).Value = "Occupied"
$ws10.Cells.Item(77, 1).Value = "2026-02-01"
$ws10.Cells.Item(77, 2).Value = "20:18:41"
$ws10.Cells.Item(77, 3).Value = "20:00"
$ws10.Cells.Item(77, 4).Value = "Bedroom"
$ws10.Cells.Item(77, 5).Value = "In Bed"
$ws10.Cells.Item(77, 6).Value = "Occupied"
$ws10.Cells.Item(78, 1).Value = "2026-02-01"
$ws10.Cells.Item(78, 2).Value = "20:18:43"
$ws10.Cells.Item(78, 3).Value = "20:00"
$ws10.Cells.Item(78, 4).Value = "Bedroom"
$ws10.Cells.Item(78, 5).Value = "In Bed"
$ws10.Cells.Item(78, 6).Value = "Occupied"
$ws10.Cells.Item(79, 1).Value = "2026-02-01"
$ws10.Cells.Item(79, 2).Value = "20:18:45"
$ws10.Cells.Item(79, 3).Value = "20:00"
$ws10.Cells.Item(79, 4).Value = "Bedroom"
$ws10.Cells.Item(79, 5).Value = "In Bed"
$ws10.Cells.Item(79, 6).Value = "Occupied"
$ws10.Cells.Item(80, 1).Value = "2026-02-01"
$ws10.Cells.Item(80, 2).Value = "20:18:46"
$ws10.Cells.Item(80, 3).Value = "20:00"
$ws10.Cells.Item(80, 4).Value = "Bedroom"
$ws10.Cells.Item(80, 5).Value = "In Bed"
$ws10.Cells.Item(80, 6).Value = "Occupied"
$ws10.Cells.Item(81, 1).Value = "2026-02-01"
$ws10.Cells.Item(81, 2).Value = "20:18:48"
$ws10.Cells.Item(81, 3).Value = "20:00"
$ws10.Cells.Item(81, 4).Value = "Bedroom"
$ws10.Cells.Item(81, 5).Value = "In Bed"
$ws10.Cells.Item(81, 6).Value = "Occupied"
$ws10.Cells.Item(82, 1).Value = "2026-02-01"
$ws10.Cells.Item(82, 2).Value = "20:18:57"
$ws10.Cells.Item(82, 3).Value = "20:00"
$ws10.Cells.Item(82, 4).Value = "Bedroom"
$ws10.Cells.Item(82, 5).Value = "In Bed"
$ws10.Cells.Item(82, 6).Value = "Occupied"
$ws10.Cells.Item(83, 1).Value = "2026-02-01"
$ws10.Cells.Item(83, 2).Value = "20:19:11"
$ws10.Cells.Item(83, 3).Value = "20:00"
$ws10.Cells.Item(83, 4).Value = "Bedroom"
$ws10.Cells.Item(83, 5).Value = "In Bed"
$ws10.Cells.Item(83, 6).Value = "Occupied"
$ws10.Cells.Item(84, 1).Value = "2026-02-01"
$ws10.Cells.Item(84, 2).Value = "20:19:19"
$ws10.Cells.Item(84, 3).Value = "20:00"
$ws10.Cells.Item(84, 4).Value = "Bedroom"
$ws10.Cells.Item(84, 5).Value = "In Bed"
$ws10.Cells.Item(84, 6).Value = "Occupied"
$ws10.Cells.Item(85, 1).Value = "2026-02-01"
$ws10.Cells.Item(85, 2).Value = "20:19:24"
$ws10.Cells.Item(85, 3).Value = "20:00"
$ws10.Cells.Item(85, 4).Value = "Bedroom"
$ws10.Cells.Item(85, 5).Value = "In Bed"
$ws10.Cells.Item(85, 6).Value = "Occupied"
$ws10.Cells.Item(86, 1).Value = "2026-02-01"
$ws10.Cells.Item(86, 2).Value = "20:19:27"
$ws10.Cells.Item(86, 3).Value = "20:00"
$ws10.Cells.Item(86, 4).Value = "Bedroom"
$ws10.Cells.Item(86, 5).Value = "In Bed"
$ws10.Cells.Item(86, 6).Value = "Occupied"
$ws10.Cells.Item(87, 1).Value = "2026-02-01"
$ws10.Cells.Item(87, 2).Value = "20:19:32"
$ws10.Cells.Item(87, 3).Value = "20:00"
$ws10.Cells.Item(87, 4).Value = "Bedroom"
$ws10.Cells.Item(87, 5).Value = "In Bed"
$ws10.Cells.Item(87, 6).Value = "Occupied"
$ws10.Cells.Item(88, 1).Value = "2026-02-01"
$ws10.Cells.Item(88, 2).Value = "20:19:34"
$ws10.Cells.Item(88, 3).Value = "20:00"
$ws10.Cells.Item(88, 4).Value = "Bedroom"
$ws10.Cells.Item(88, 5).Value = "In Bed"
$ws10.Cells.Item(88, 6).Value = "Occupied"

# --- Sheet: mmWave(BR) -> add rows 72-84 ---
$ws8 = $wb.Worksheets.Item("mmWave(BR)")
$ws8.Range("A72:A84").NumberFormat = "@"
$ws8.Cells.Item(72, 1).Value = "2026-02-01"
$ws8.Cells.Item(72, 2).Value = "20:18:36"
$ws8.Cells.Item(72, 3).Value = "20:00"
$ws8.Cells.Item(72, 4).Value = "Bedroom"
$ws8.Cells.Item(72, 5).Value = 1
$ws8.Cells.Item(72, 6).Value = "Occupied"
$ws8.Cells.Item(73, 1).Value = "2026-02-01"
$ws8.Cells.Item(73, 2).Value = "20:18:42"
$ws8.Cells.Item(73, 3).Value = "20:00"
$ws8.Cells.Item(73, 4).Value = "Bedroom"
$ws8.Cells.Item(73, 5).Value = 2
$ws8.Cells.Item(73, 6).Value = "Occupied"
$ws8.Cells.Item(74, 1).Value = "2026-02-01"
$ws8.Cells.Item(74, 2).Value = "20:18:44"
$ws8.Cells.Item(74, 3).Value = "20:00"
$ws8.Cells.Item(74, 4).Value = "Bedroom"
$ws8.Cells.Item(74, 5).Value = 8
$ws8.Cells.Item(74, 6).Value = "Occupied"
$ws8.Cells.Item(75, 1).Value = "2026-02-01"
$ws8.Cells.Item(75, 2).Value = "20:18:46"
$ws8.Cells.Item(75, 3).Value = "20:00"
$ws8.Cells.Item(75, 4).Value = "Bedroom"
$ws8.Cells.Item(75, 5).Value = 2
$ws8.Cells.Item(75, 6).Value = "Occupied"
$ws8.Cells.Item(76, 1).Value = "2026-02-01"
$ws8.Cells.Item(76, 2).Value = "20:18:47"
$ws8.Cells.Item(76, 3).Value = "20:00"
$ws8.Cells.Item(76, 4).Value = "Bedroom"
$ws8.Cells.Item(76, 5).Value = 3
$ws8.Cells.Item(76, 6).Value = "Occupied"
$ws8.Cells.Item(77, 1).Value = "2026-02-01"
$ws8.Cells.Item(77, 2).Value = "20:18:49"
$ws8.Cells.Item(77, 3).Value = "20:00"
$ws8.Cells.Item(77, 4).Value = "Bedroom"
$ws8.Cells.Item(77, 5).Value = 2
$ws8.Cells.Item(77, 6).Value = "Occupied"
$ws8.Cells.Item(78, 1).Value = "2026-02-01"
$ws8.Cells.Item(78, 2).Value = "20:18:58"
$ws8.Cells.Item(78, 3).Value = "20:00"
$ws8.Cells.Item(78, 4).Value = "Bedroom"
$ws8.Cells.Item(78, 5).Value = 1
$ws8.Cells.Item(78, 6).Value = "Occupied"
$ws8.Cells.Item(79, 1).Value = "2026-02-01"
$ws8.Cells.Item(79, 2).Value = "20:19:13"
$ws8.Cells.Item(79, 3).Value = "20:00"
$ws8.Cells.Item(79, 4).Value = "Bedroom"
$ws8.Cells.Item(79, 5).Value = 2
$ws8.Cells.Item(79, 6).Value = "Occupied"
$ws8.Cells.Item(80, 1).Value = "2026-02-01"
$ws8.Cells.Item(80, 2).Value = "20:19:20"
$ws8.Cells.Item(80, 3).Value = "20:00"
$ws8.Cells.Item(80, 4).Value = "Bedroom"
$ws8.Cells.Item(80, 5).Value = 1
$ws8.Cells.Item(80, 6).Value = "Occupied"
$ws8.Cells.Item(81, 1).Value = "2026-02-01"
$ws8.Cells.Item(81, 2).Value = "20:19:25"
$ws8.Cells.Item(81, 3).Value = "20:00"
$ws8.Cells.Item(81, 4).Value = "Bedroom"
$ws8.Cells.Item(81, 5).Value = 2
$ws8.Cells.Item(81, 6).Value = "Occupied"
$ws8.Cells.Item(82, 1).Value = "2026-02-01"
$ws8.Cells.Item(82, 2).Value = "20:19:28"
$ws8.Cells.Item(82, 3).Value = "20:00"
$ws8.Cells.Item(82, 4).Value = "Bedroom"
$ws8.Cells.Item(82, 5).Value = 1
$ws8.Cells.Item(82, 6).Value = "Occupied"
$ws8.Cells.Item(83, 1).Value = "2026-02-01"
$ws8.Cells.Item(83, 2).Value = "20:19:33"
$ws8.Cells.Item(83, 3).Value = "20:00"
$ws8.Cells.Item(83, 4).Value = "Bedroom"
$ws8.Cells.Item(83, 5).Value = 2
$ws8.Cells.Item(83, 6).Value = "Occupied"
$ws8.Cells.Item(84, 1).Value = "2026-02-01"
$ws8.Cells.Item(84, 2).Value = "20:19:35"
$ws8.Cells.Item(84, 3).Value = "20:00"
$ws8.Cells.Item(84, 4).Value = "Bedroom"
$ws8.Cells.Item(84, 5).Value = 3
$ws8.Cells.Item(84, 6).Value = "Occupied"

# --- Sheet: mmWave(HR) -> add rows 72-84 ---
$ws9 = $wb.Worksheets.Item("mmWave(HR)")
$ws9.Range("A72:A84").NumberFormat = "@"
$ws9.Cells.Item(72, 1).Value = "2026-02-01"
$ws9.Cells.Item(72, 2).Value = "20:18:36"
$ws9.Cells.Item(72, 3).Value = "20:00"
$ws9.Cells.Item(72, 4).Value = "Bedroom"
$ws9.Cells.Item(72, 5).Value = 49
$ws9.Cells.Item(72, 6).Value = "Occupied"
$ws9.Cells.Item(73, 1).Value = "2026-02-01"
$ws9.Cells.Item(73, 2).Value = "20:18:42"
$ws9.Cells.Item(73, 3).Value = "20:00"
$ws9.Cells.Item(73, 4).Value = "Bedroom"
$ws9.Cells.Item(73, 5).Value = 50
$ws9.Cells.Item(73, 6).Value = "Occupied"
$ws9.Cells.Item(74, 1).Value = "2026-02-01"
$ws9.Cells.Item(74, 2).Value = "20:18:43"
$ws9.Cells.Item(74, 3).Value = "20:00"
$ws9.Cells.Item(74, 4).Value = "Bedroom"
$ws9.Cells.Item(74, 5).Value = 56
$ws9.Cells.Item(74, 6).Value = "Occupied"
$ws9.Cells.Item(75, 1).Value = "2026-02-01"
$ws9.Cells.Item(75, 2).Value = "20:18:45"
$ws9.Cells.Item(75, 3).Value = "20:00"
$ws9.Cells.Item(75, 4).Value = "Bedroom"
$ws9.Cells.Item(75, 5).Value = 50
$ws9.Cells.Item(75, 6).Value = "Occupied"
$ws9.Cells.Item(76, 1).Value = "2026-02-01"
$ws9.Cells.Item(76, 2).Value = "20:18:47"
$ws9.Cells.Item(76, 3).Value = "20:00"
$ws9.Cells.Item(76, 4).Value = "Bedroom"
$ws9.Cells.Item(76, 5).Value = 51
$ws9.Cells.Item(76, 6).Value = "Occupied"
$ws9.Cells.Item(77, 1).Value = "2026-02-01"
$ws9.Cells.Item(77, 2).Value = "20:18:48"
$ws9.Cells.Item(77, 3).Value = "20:00"
$ws9.Cells.Item(77, 4).Value = "Bedroom"
$ws9.Cells.Item(77, 5).Value = 50
$ws9.Cells.Item(77, 6).Value = "Occupied"
$ws9.Cells.Item(78, 1).Value = "2026-02-01"
$ws9.Cells.Item(78, 2).Value = "20:18:57"
$ws9.Cells.Item(78, 3).Value = "20:00"
$ws9.Cells.Item(78, 4).Value = "Bedroom"
$ws9.Cells.Item(78, 5).Value = 49
$ws9.Cells.Item(78, 6).Value = "Occupied"
$ws9.Cells.Item(79, 1).Value = "2026-02-01"
$ws9.Cells.Item(79, 2).Value = "20:19:12"
$ws9.Cells.Item(79, 3).Value = "20:00"
$ws9.Cells.Item(79, 4).Value = "Bedroom"
$ws9.Cells.Item(79, 5).Value = 50
$ws9.Cells.Item(79, 6).Value = "Occupied"
$ws9.Cells.Item(80, 1).Value = "2026-02-01"
$ws9.Cells.Item(80, 2).Value = "20:19:20"
$ws9.Cells.Item(80, 3).Value = "20:00"
$ws9.Cells.Item(80, 4).Value = "Bedroom"
$ws9.Cells.Item(80, 5).Value = 49
$ws9.Cells.Item(80, 6).Value = "Occupied"
$ws9.Cells.Item(81, 1).Value = "2026-02-01"
$ws9.Cells.Item(81, 2).Value = "20:19:25"
$ws9.Cells.Item(81, 3).Value = "20:00"
$ws9.Cells.Item(81, 4).Value = "Bedroom"
$ws9.Cells.Item(81, 5).Value = 50
$ws9.Cells.Item(81, 6).Value = "Occupied"
$ws9.Cells.Item(82, 1).Value = "2026-02-01"
$ws9.Cells.Item(82, 2).Value = "20:19:27"
$ws9.Cells.Item(82, 3).Value = "20:00"
$ws9.Cells.Item(82, 4).Value = "Bedroom"
$ws9.Cells.Item(82, 5).Value = 49
$ws9.Cells.Item(82, 6).Value = "Occupied"
$ws9.Cells.Item(83, 1).Value = "2026-02-01"
$ws9.Cells.Item(83, 2).Value = "20:19:32"
$ws9.Cells.Item(83, 3).Value = "20:00"
$ws9.Cells.Item(83, 4).Value = "Bedroom"
$ws9.Cells.Item(83, 5).Value = 50
$ws9.Cells.Item(83, 6).Value = "Occupied"
$ws9.Cells.Item(84, 1).Value = "2026-02-01"
$ws9.Cells.Item(84, 2).Value = "20:19:35"
$ws9.Cells.Item(84, 3).Value = "20:00"
$ws9.Cells.Item(84, 4).Value = "Bedroom"
$ws9.Cells.Item(84, 5).Value = 51
$ws9.Cells.Item(84, 6).Value = "Occupied"
